$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.345.67"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "2.609.11"
$ws.Range("E3").Value = "  +2.02%  "

$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "589.60"
$ws.Range("E5").Value = "  +3.80%  "

$ws.Range("D6").Value = "149.57"
$ws.Range("E6").Value = "  +1.99%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  +2.09%  "

$ws.Range("E9").Value = "  +4.50%  "

$ws.Range("D10").Value = "5.68"
$ws.Range("E10").Value = "  +1.88%  "

$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("E12").Value = "  +1.70%  "

$ws.Range("D13").Value = "27.88"
$ws.Range("E13").Value = "  +1.27%  "

$ws.Range("D14").Value = "3.072.39"
$ws.Range("E14").Value = "  +2.00%  "

$ws.Range("D15").Value = "63.330.89"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").Value = "0.0000150"
$ws.Range("E16").Value = "  +4.38%  "

$ws.Range("D17").Value = "2.591.20"
$ws.Range("E17").Value = "  +1.21%  "

$ws.Range("D18").Value = "11.49"
$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("D19").Value = "346.22"
$ws.Range("E19").Value = "  +3.18%  "

$ws.Range("D20").Value = "4.47"
$ws.Range("E20").Value = "  +2.99%  "

$ws.Range("D21").Value = "6.91"
$ws.Range("E21").Value = "  +1.58%  "

$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("E23").Value = "  -3.39%  "

$ws.Range("D24").Value = "66.93"
$ws.Range("E24").Value = "  +2.56%  "

$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "2.676.53"
$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.171"
$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("D27").Value = "1.64"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").Value = "8.25"
$ws.Range("E28").Value = "  +12.33%  "

$ws.Range("D29").Value = "8.56"
$ws.Range("E29").Value = "  +0.93%  "

$ws.Range("E30").Value = "  +1.67%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  +8.41%  "

$ws.Range("D33").Value = "0.0₃0839"
$ws.Range("E33").Value = "  +2.39%  "

$ws.Range("D34").Value = "469.39"
$ws.Range("E34").Value = "  +15.27%  "

$ws.Range("D35").Value = "1.65"
$ws.Range("E35").Value = "  +4.89%  "

$ws.Range("D36").Value = "177.23"
$ws.Range("E36").Value = "  +0.75%  "

$ws.Range("D37").Value = "0.409"
$ws.Range("E37").Value = "  +2.66%  "

$ws.Range("D38").Value = "19.38"
$ws.Range("E38").Value = "  +2.08%  "

$ws.Range("D39").Value = "4.68"
$ws.Range("E39").Value = "  +6.76%  "

$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("E41").Value = "  +0.96%  "

$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D43").Value = "153.62"
$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("D44").Value = "3.87"
$ws.Range("E44").Value = "  +2.36%  "

$ws.Range("D45").Value = "21.24"
$ws.Range("E45").Value = "  +1.00%  "

$ws.Range("D46").Value = "0.0557"
$ws.Range("E46").Value = "  +5.83%  "

$ws.Range("D47").Value = "0.619"
$ws.Range("E47").Value = "  +1.94%  "

$ws.Range("D48").Value = "0.0980"
$ws.Range("E48").Value = "  +1.89%  "

$ws.Range("D49").Value = "0.0243"
$ws.Range("E49").Value = "  +1.73%  "

$ws.Range("D50").Value = "1.78"
$ws.Range("E50").Value = "  +0.12%  "

$ws.Range("E51").Value = "  +0.83%  "
